# Auto-generated edit script: updates the cryptos price/volume table
# (commit: "Updated cryptos list on Sat Jun 29 10:26:42 UTC 2024 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must stay TEXT (Price/Volume columns hold
# strings such as "60.909.00" or "1.00" - plain assignment would let Excel
# coerce them to numbers and silently drop meaningful trailing/format info).
# Prefixing with an apostrophe forces text entry; resetting the Style back
# to "Normal" afterwards avoids leaving a stray text-number-format style
# behind (keeps the cell style identical to the original General style).
function Set-TextValue($range, $text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") '60.909.00'
Set-TextValue $ws.Range("E2") '  -0.94%  '
# Row 3
Set-TextValue $ws.Range("D3") '3.393.01'
Set-TextValue $ws.Range("E3") '  -1.51%  '
# Row 4
Set-TextValue $ws.Range("D4") '0.999'
Set-TextValue $ws.Range("E4") '  +0.03%  '
# Row 5
Set-TextValue $ws.Range("D5") '571.85'
Set-TextValue $ws.Range("E5") '  -0.92%  '
# Row 6
Set-TextValue $ws.Range("D6") '142.65'
Set-TextValue $ws.Range("E6") '  -1.77%  '
# Row 7
Set-TextValue $ws.Range("D7") '3.393.28'
Set-TextValue $ws.Range("E7") '  -1.48%  '
# Row 8
Set-TextValue $ws.Range("E8") '  +0.03%  '
# Row 9
Set-TextValue $ws.Range("E9") '  -0.54%  '
# Row 10
Set-TextValue $ws.Range("E10") '  -1.78%  '
# Row 11
Set-TextValue $ws.Range("E11") '  -2.71%  '
# Row 12
Set-TextValue $ws.Range("E12") '  +1.47%  '
# Row 13
Set-TextValue $ws.Range("D13") '3.972.32'
Set-TextValue $ws.Range("E13") '  -1.37%  '
# Row 14
Set-TextValue $ws.Range("E14") '  +2.02%  '
# Row 15
Set-TextValue $ws.Range("D15") '28.19'
Set-TextValue $ws.Range("E15") '  -1.32%  '
# Row 16
Set-TextValue $ws.Range("E16") '  -1.53%  '
# Row 17
Set-TextValue $ws.Range("D17") '3.393.77'
Set-TextValue $ws.Range("E17") '  -1.54%  '
# Row 18
Set-TextValue $ws.Range("D18") '60.988.53'
Set-TextValue $ws.Range("E18") '  -1.01%  '
# Row 19
Set-TextValue $ws.Range("D19") '6.15'
Set-TextValue $ws.Range("E19") '  -2.88%  '
# Row 20
Set-TextValue $ws.Range("D20") '13.84'
Set-TextValue $ws.Range("E20") '  -3.52%  '
# Row 21
Set-TextValue $ws.Range("D21") '8.98'
Set-TextValue $ws.Range("E21") '  -5.07%  '
# Row 22
Set-TextValue $ws.Range("D22") '382.87'
Set-TextValue $ws.Range("E22") '  -4.08%  '
# Row 23
Set-TextValue $ws.Range("D23") '0.559'
Set-TextValue $ws.Range("E23") '  -1.84%  '
# Row 24
Set-TextValue $ws.Range("D24") '74.25'
Set-TextValue $ws.Range("E24") '  +0.31%  '
# Row 25
Set-TextValue $ws.Range("D25") '1.00'
Set-TextValue $ws.Range("E25") '  -0.01%  '
# Row 26
Set-TextValue $ws.Range("E26") '  -4.40%  '
# Row 27
Set-TextValue $ws.Range("D27") '3.528.52'
Set-TextValue $ws.Range("E27") '  -1.58%  '
# Row 28
Set-TextValue $ws.Range("E28") '  -0.93%  '
# Row 29
Set-TextValue $ws.Range("E29") '  -0.11%  '
# Row 30
Set-TextValue $ws.Range("E30") '  -3.13%  '
# Row 31
Set-TextValue $ws.Range("E31") '  -3.66%  '
# Row 32
Set-TextValue $ws.Range("E32") '  -1.82%  '
# Row 33
Set-TextValue $ws.Range("D33") '1.41'
Set-TextValue $ws.Range("E33") '  -3.00%  '
# Row 34
Set-TextValue $ws.Range("E34") '  -0.04%  '
# Row 35
Set-TextValue $ws.Range("D35") '23.51'
Set-TextValue $ws.Range("E35") '  -1.76%  '
# Row 36
Set-TextValue $ws.Range("D36") '6.99'
Set-TextValue $ws.Range("E36") '  -0.65%  '
# Row 37
Set-TextValue $ws.Range("D37") '167.79'
Set-TextValue $ws.Range("E37") '  +0.40%  '
# Row 38
Set-TextValue $ws.Range("D38") '3.423.67'
Set-TextValue $ws.Range("E38") '  -1.35%  '
# Row 39
Set-TextValue $ws.Range("E39") '  -3.17%  '
# Row 40
Set-TextValue $ws.Range("E40") '  -4.72%  '
# Row 41
$ws.Range("B41").Value = 'EnergySwap'
$ws.Range("C41").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range("D41") '27.81'
Set-TextValue $ws.Range("E41") '  +1.95%  '
# Row 42
$ws.Range("B42").Value = 'Hedera'
$ws.Range("C42").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range("D42") '0.0773'
Set-TextValue $ws.Range("E42") '  -2.39%  '
# Row 43
Set-TextValue $ws.Range("E43") '  -2.81%  '
# Row 44
Set-TextValue $ws.Range("E44") '  +0.05%  '
# Row 45
$ws.Range("B45").Value = 'Filecoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range("D45") '4.44'
Set-TextValue $ws.Range("E45") '  -1.82%  '
# Row 46
$ws.Range("B46").Value = 'Stacks'
$ws.Range("C46").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range("D46") '1.67'
Set-TextValue $ws.Range("E46") '  -3.62%  '
# Row 47
$ws.Range("B47").Value = 'ONDO'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
Set-TextValue $ws.Range("D47") '1.14'
Set-TextValue $ws.Range("E47") '  -1.22%  '
# Row 48
$ws.Range("B48").Value = 'Maker'
$ws.Range("C48").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws.Range("D48") '2.476.20'
Set-TextValue $ws.Range("E48") '  -5.10%  '
# Row 49
$ws.Range("B49").Value = 'Cosmos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue $ws.Range("D49") '6.82'
Set-TextValue $ws.Range("E49") '  -1.91%  '
# Row 50
$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range("D50") '22.98'
Set-TextValue $ws.Range("E50") '  -0.60%  '
# Row 51
$ws.Range("B51").Value = 'VeChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range("D51") '0.0267'
Set-TextValue $ws.Range("E51") '  +1.18%  '
